$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume figures for the latest symbol-list refresh.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '306.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.15%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.88'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.08%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.058'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.14%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08090'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.94%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.951'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.82%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.156'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.49%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '7.768'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.47%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9285'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.29%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1366'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.83%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1905'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.65%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09225'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.62%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03534'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '3.33%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09853'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.18%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001431'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.88%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005816'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.01%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.565'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.55%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.976'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.32%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3444'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.21%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1345'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.37%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.889'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-3.12%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2512'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.58%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04420'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.57%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001222'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.59%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004776'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.70%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '31.95%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003124'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '4.06%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01949'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.85%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04993'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '5.49%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01097'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '14.06%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007630'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4.03%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1381'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.39%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002100'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.56%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01080'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.76%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006377'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.35%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.09%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '65.22'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '1.15%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001189'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-20.19%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.09%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.09%'
